# edit.ps1
# Applies the commit "Added a few more slots" to celestial-king document:
#  1. Insert a new "Meta description" paragraph right after the title (Heading1) paragraph.
#  2. Remove the duplicated bold title paragraph near the end of the document.
#  3. Replace the italic "meta description" paragraph's text at the end with the new
#     "Prompt: ..." image-generation prompt text (keeping italics).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 0: Replace the italic "meta description" paragraph text at the bottom
# with the new image-generation prompt text (formatting/italics preserved).
# Done first, while the matched phrase is still unique in the document.
#
# We avoid Find.Execute's ReplaceWith here because it runs the replacement
# text through AutoCorrect/AutoFormat (turning straight apostrophes into
# curly ones); instead we insert the new text manually and then delete the
# old text, re-applying italics explicitly afterwards.
# ---------------------------------------------------------------------------

$promptText = "Prompt: Please create a feature image for Celestial King online slot game that is cartoon-style and features a happy Maya warrior with glasses. Our team of talented graphic designers can create a stunning feature image for Celestial King online slot game. We understand the need for an eye-catching design that perfectly captures the essence of the game. Based on your prompt, we will create an image that features a happy Maya warrior with glasses in a cartoon-style. The warrior will be depicted in a vibrant color palette that matches the game's ambiance. The lively expression on the warrior's face is sure to capture the attention of potential players. The image will showcase the Celestial King holding a scepter with his adorable green eyes shining with amusement. The background will be steeped in vibrant purple and golden clouds with an Asian theme that perfectly matches this unique slot game's setting. Our team can create an incredible feature image that will effectively promote Celestial King online slot game, attracting potential players."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Explore the Asian-inspired Celestial King slot game with stunning graphics*") {
        $oldLen = $p.Range.Text.Length - 1   # exclude the trailing paragraph mark
        $insPoint = $d.Range($p.Range.Start, $p.Range.Start)
        $insPoint.InsertAfter($promptText)
        $deleteRange = $d.Range($p.Range.Start + $promptText.Length, $p.Range.Start + $promptText.Length + $oldLen)
        $deleteRange.Delete()
        $newRange = $d.Range($p.Range.Start, $p.Range.Start + $promptText.Length)
        $newRange.Font.Italic = $true
        break
    }
}

# ---------------------------------------------------------------------------
# Step 1: Build the new "Meta description" paragraph.
# We build it next to an existing plain (Normal-style) paragraph so that it
# naturally picks up plain/Normal paragraph formatting (no heading style),
# then move it into place right after the title paragraph.
# ---------------------------------------------------------------------------

$boldLabel = "Meta description"
$metaBody  = ": Explore the Asian-inspired Celestial King slot game with stunning graphics and jackpots. Read our review and play free at top online casinos."

# Paragraph 3 ("Looking for a slot game...") is a normal body paragraph - insert
# a fresh paragraph right after it so the new paragraph inherits plain formatting.
$anchorPara = $d.Paragraphs.Item(3)
$anchorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(4)

# Fill in the text (label + body) in the new, still-empty paragraph.
$insertPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$insertPoint.InsertBefore($boldLabel + $metaBody)

$newPara = $d.Paragraphs.Item(4)

# Make just the "Meta description" label bold.
$boldRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $boldLabel.Length)
$boldRange.Font.Bold = $true

# Move the whole new paragraph (including its paragraph mark) so it sits right
# after the title paragraph (paragraph 1).
$newPara = $d.Paragraphs.Item(4)
$fullNewRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$fullNewRange.Cut()

$titlePara = $d.Paragraphs.Item(1)
$target = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$target.Paste()

# ---------------------------------------------------------------------------
# Step 2: Remove the duplicated bold title paragraph near the end of the doc
# ("Play Celestial King Free - Pros and Cons of Online Slot Game").
# ---------------------------------------------------------------------------

# Locate the duplicate paragraph directly: it's the paragraph whose full text
# equals the bold title (the first such paragraph is the real Heading1 title,
# so we look for the later, non-heading occurrence).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Celestial King Free - Pros and Cons of Online Slot Game`r" -or $p.Range.Text -eq "Play Celestial King Free - Pros and Cons of Online Slot Game") {
        if ($p.Style.NameLocal -ne "Heading 1") {
            $p.Range.Delete()
            break
        }
    }
}
